# Circle Language Spec Plan / Clarify Command as a Concept - Project Summary
#
# 1) Move the "_GoBack" bookmark (Word's "last edit position" marker) from
#    in front of the "Goal" heading to inside "Project Summary" - splitting
#    that run into "Project Sum" | bookmark | "mary", exactly where the
#    author's cursor was left after the last edit.
# 2) Bump the Heading2 style's font from Arial 16pt to Calibri 18pt.

$d = $word.ActiveDocument

# --- Remove the old "_GoBack" bookmark (currently sits before "Goal") ----
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Re-create "_GoBack" right after "Project Sum", splitting the run ----
$rng = $d.Content
$found = $rng.Find.Execute("Project Sum", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)   # wdCollapseEnd -> zero-length range right after "Project Sum"
    $d.Bookmarks.Add("_GoBack", $rng)
}

# --- Heading2 style: Arial -> Calibri (ascii/hAnsi), 16pt -> 18pt --------
$style = $d.Styles("Heading 2")
$style.Font.Name = "Calibri"
$style.Font.Size = 18
